$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd() -eq $needle) {
            return $p
        }
    }
    return $null
}

# --- Change 1: insert two new member-function declarations
#     ("getCabezaEstaciones" and "getSiguienteLinea") right before the
#     existing "    int getNumEstaciones() const;" declaration. ---
$target1 = Find-ParagraphByText $d "    int getNumEstaciones() const;"
$target1.Range.InsertBefore("    Estacion* getCabezaEstaciones() const;`r    Linea* getSiguienteLinea() const;`r")

# --- Change 2: insert a new blank paragraph right after the
#     "#endif" line (before the existing trailing blank paragraph). ---
$target2 = Find-ParagraphByText $d "#endif"
$target2.Range.InsertParagraphAfter()
